$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert the new columns -------------------------------------------------
# Two new columns right after the existing "capacity_plant_kw" column (C),
# shifting the old D..M block to F..N.
$ws.Range("D1").EntireColumn.Insert()
$ws.Range("D1").EntireColumn.Insert()

# One more new column right before the old "capacity_%" column (now at N),
# shifting it (and generation_%) one slot to the right.
$ws.Range("N1").EntireColumn.Insert()

# --- Header row --------------------------------------------------------------
$ws.Range("C1").Value = "plant_capacity_power_kw"
$ws.Range("D1").Value = "plant_capacity_force_kwh"
$ws.Range("E1").Value = "plant_generate_force_kwh"
$ws.Range("N1").Value = "load_factor"

# --- New formula columns, rows 2-8 ------------------------------------------
for ($r = 2; $r -le 8; $r++) {
    $ws.Range("D$r").Formula = "=C$r* 365 * 24"
    $ws.Range("E$r").Formula = "=D$r*N$r"
    $ws.Range("N$r").Formula = "=M$r/K$r"
}
